$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4351088404655457
$ws.Range("B1").Value = 1.026312232017517
$ws.Range("C1").Value = 4.622306823730469
$ws.Range("D1").Value = 1.442337155342102
$ws.Range("E1").Value = 1.17347526550293
